# Apply updated values to the "Xr_Results" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ---
$ws.Range("C7").Value = 151
$ws.Range("D7").Value = 4401
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 4401
$ws.Range("K7").Value = 170
$ws.Range("L7").Value = 4880

# --- Row 8 ---
$ws.Range("C8").Value = 66
$ws.Range("D8").Value = 4233
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 4306
$ws.Range("K8").Value = 82
$ws.Range("L8").Value = 4864

# --- Row 11 ---
$ws.Range("C11").Value = 0.08988095238095238
$ws.Range("G11").Value = 0.0
$ws.Range("K11").Value = 0.08374384236453201

# --- Row 12 ---
$ws.Range("C12").Value = 0.039285714285714285
$ws.Range("G12").Value = 0.0
$ws.Range("K12").Value = 0.04039408866995074

# --- Row 18 ---
$ws.Range("C18").Value = 106
$ws.Range("D18").Value = 4995
$ws.Range("G18").Value = 76
$ws.Range("H18").Value = 4620
$ws.Range("K18").Value = 503
$ws.Range("L18").Value = 23297

# --- Row 19 ---
$ws.Range("C19").Value = 76
$ws.Range("D19").Value = 4996
$ws.Range("G19").Value = 69
$ws.Range("H19").Value = 4371
$ws.Range("K19").Value = 293
$ws.Range("L19").Value = 22770

# --- Row 22 ---
$ws.Range("C22").Value = 0.061378112333526344
$ws.Range("G22").Value = 0.2585034013605442
$ws.Range("K22").Value = 0.08776827778747165

# --- Row 23 ---
$ws.Range("C23").Value = 0.04400694846554719
$ws.Range("G23").Value = 0.23469387755102042
$ws.Range("K23").Value = 0.0511254580352469
